# Fix bug: "procurar nome produto" p -> h3 (commit message).
# Adds a GoogleTrends column (K) and refreshes the scraped data for the
# first two products while zeroing-out / marking "NA" the rows that no
# longer have fresh data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column K: header ---------------------------------------------
# Copy formatting (bold, border, centered) from an existing header cell
# so the new header matches the rest of row 1.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K1").Value = "GoogleTrends"

# --- Row 2: bbs airsoft --------------------------------------------------
$ws.Range("D2").Value = "2022-05-20 19:23:22"
$ws.Range("E2").Value = "17.113 resultado"
$ws.Range("G2").Value = 85.33333333333333
$ws.Range("H2").Value = 89
$ws.Range("I2").Value = 3381
$ws.Range("J2").Value = 4454
$ws.Range("K2").Value = "https://trends.google.com.br/trends/explore?geo=BR&q=bbs airsoft"

# --- Row 3: alvo tiro -----------------------------------------------------
$ws.Range("D3").Value = "2022-05-20 19:23:22"
$ws.Range("E3").Value = "10.822 resultado"
$ws.Range("G3").Value = 35
$ws.Range("H3").Value = 28
$ws.Range("I3").Value = 181
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = "https://trends.google.com.br/trends/explore?geo=BR&q=alvo tiro"

# --- Rows 4-11: no fresh scrape data this run -> zeroed out, marked NA ---
for ($r = 4; $r -le 11; $r++) {
    $ws.Range("D$r").Value = "2022-05-20 19:23:22"
    $ws.Range("E$r").Value = 0
    $ws.Range("F$r").Value = 0
    $ws.Range("G$r").Value = 0
    $ws.Range("H$r").Value = 0
    $ws.Range("I$r").Value = 0
    $ws.Range("J$r").Value = 0
    $ws.Range("K$r").Value = "NA"
}
